$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "Globo"
$ws.Range("B41").Value = "RJ TV 1"
$ws.Range("C41").Value = "Defesa Civil"
$ws.Range("D41").Value = "2025-04-04T12:09"
$ws.Range("E41").Value = "Positivo"
$ws.Range("F41").Value = "Preparação para a frente fria. Cidades do Norte Fluminense se mobilizam diante da previsão de chuva intensa. Repórter *ao vivo* do Cais da Lapa. Equipes da Defesa Civil em alerta. Previsão de chuvas com volume considerável e ventos fortes. Para prevenir, Prefeitura começou a fazer ações preventivas como limpeza de bueiros e de canais. Alerta se mantém para possível queda de árvores e pontos de alagamento. Entrevista com coord. Centro de Monitoramento de Desastres, Leandro Freitas, e com subsecretário da Defesa Civil de Campos, major Edson Pessanha. "

$ws.Range("A42").Value = "Globo"
$ws.Range("B42").Value = "RJ TV 1"
$ws.Range("C42").Value = "Defesa Civil"
$ws.Range("D42").Value = "2025-04-04T12:26"
$ws.Range("E42").Value = "Neutro"
$ws.Range("F42").Value = "Escombros de prédio histórico serão retirados. Operação vai interditar rua importante no Centro de Campos durante o final de semana. Repórter *ao vivo* do local. Prédio foi interditado pela Defesa Civil, após telhado ceder. Prédio foi demolido com uso de um caminhão específico após horário comercial para não atrapalhar as vendas. Vídeo com depoimento do subsecretário da Defesa Civil, major Edson Pessanha. "
